$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.205.49'
$ws.Range("E2").Value = '  +1.36%  '
$ws.Range("D3").Value = '1.651.60'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  -0.98%  '
$ws.Range("D5").Value = "'219.98"
$ws.Range("E5").Value = '  +1.46%  '
$ws.Range("D6").Value = "'0.501"
$ws.Range("E6").Value = '  -0.34%  '
$ws.Range("E7").Value = '  -0.82%  '
$ws.Range("E8").Value = '  +0.92%  '
$ws.Range("E9").Value = '  -0.20%  '
$ws.Range("D10").Value = "'19.69"
$ws.Range("E10").Value = '  +2.24%  '
$ws.Range("D11").Value = "'0.0848"
$ws.Range("E11").Value = '  +0.44%  '
$ws.Range("D12").Value = '1.882.00'
$ws.Range("E12").Value = '  +0.19%  '
$ws.Range("D13").Value = '1.661.60'
$ws.Range("E13").Value = '  +0.67%  '
$ws.Range("D14").Value = "'4.21"
$ws.Range("E14").Value = '  +0.88%  '
$ws.Range("D15").Value = "'0.533"
$ws.Range("E15").Value = '  +0.68%  '
$ws.Range("D16").Value = "'66.26"
$ws.Range("E16").Value = '  +2.21%  '
$ws.Range("D17").Value = '27.167.30'
$ws.Range("E17").Value = '  +1.23%  '
$ws.Range("D18").Value = '0.0₃0738'
$ws.Range("E18").Value = '  +0.29%  '
$ws.Range("D19").Value = "'222.58"
$ws.Range("E19").Value = '  +3.45%  '
$ws.Range("E20").Value = '  -0.92%  '
$ws.Range("D21").Value = "'6.82"
$ws.Range("E21").Value = '  +8.70%  '
$ws.Range("E22").Value = '  +0.82%  '
$ws.Range("D23").Value = "'2.43"
$ws.Range("E23").Value = '  -2.61%  '
$ws.Range("E24").Value = '  -0.75%  '
$ws.Range("D25").Value = "'147.13"
$ws.Range("E25").Value = '  -0.20%  '
$ws.Range("E26").Value = '  -0.56%  '
$ws.Range("D27").Value = "'7.37"
$ws.Range("E27").Value = '  +2.04%  '
$ws.Range("E28").Value = '  +0.28%  '
$ws.Range("D29").Value = "'15.96"
$ws.Range("D30").Value = "'0.0515"
$ws.Range("E30").Value = '  +1.33%  '
$ws.Range("E31").Value = '  +1.08%  '
$ws.Range("E32").Value = '  +0.54%  '
$ws.Range("D33").Value = "'3.03"
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("E34").Value = '  +2.80%  '
$ws.Range("D35").Value = '1.266.44'
$ws.Range("E35").Value = '  -2.39%  '
$ws.Range("E36").Value = '  -0.45%  '
$ws.Range("E37").Value = '  -1.49%  '
$ws.Range("E38").Value = '  +0.43%  '
$ws.Range("E40").Value = '  -0.71%  '
$ws.Range("D41").Value = "'0.809"
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").Value = "'5.40"
$ws.Range("E42").Value = '  +1.20%  '
$ws.Range("D43").Value = '1.792.36'
$ws.Range("E43").Value = '  +0.29%  '
$ws.Range("D44").Value = "'62.10"
$ws.Range("E44").Value = '  +0.10%  '
$ws.Range("D45").Value = "'92.57"
$ws.Range("E45").Value = '  +0.69%  '
$ws.Range("D46").Value = "'2.06"
$ws.Range("E46").Value = '  -8.09%  '
$ws.Range("E47").Value = '  +0.30%  '
$ws.Range("D49").Value = "'7.65"
$ws.Range("E49").Value = '  -0.13%  '
$ws.Range("D50").Value = "'0.0977"
$ws.Range("E50").Value = '  +0.45%  '
$ws.Range("E51").Value = '  -0.58%  '

# Cells whose new values look numeric were entered with a leading
# apostrophe above so Excel keeps them as text (matching the original
# inlineStr cell type); reset style index so no stray number-format
# is left behind on those cells.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
